$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "嘉泽新能"
$ws.Range("B3").Value = "明阳智能"
$ws.Range("C3").Value = "万向钱潮"
$ws.Range("A4").Value = "立讯精密"
$ws.Range("B4").Value = "万向钱潮"
$ws.Range("C4").Value = "天赐材料"
$ws.Range("A5").Value = "上海电气"
$ws.Range("B5").Value = "中际旭创"
$ws.Range("C5").Value = "精艺股份"
$ws.Range("A6").Value = "万向钱潮"
$ws.Range("C6").Value = "和而泰"
$ws.Range("A7").Value = "山子高科"
$ws.Range("B7").Value = "天赐材料"
$ws.Range("C7").Value = "吉鑫科技"
$ws.Range("A8").Value = "吉鑫科技"
$ws.Range("B8").Value = "吉鑫科技"
$ws.Range("C8").Value = "立讯精密"
$ws.Range("B9").Value = "上海电气"
$ws.Range("A10").Value = "张江高科"
$ws.Range("B10").Value = "先导智能"
$ws.Range("C10").Value = "三花智控"
$ws.Range("A11").Value = "金风科技"
$ws.Range("B11").Value = "多氟多"
$ws.Range("A12").Value = "先导智能"
$ws.Range("B12").Value = "赛力斯"
$ws.Range("A13").Value = "多氟多"
$ws.Range("B13").Value = "嘉泽新能"
$ws.Range("C13").Value = "嘉泽新能"
$ws.Range("A14").Value = "华建集团"
$ws.Range("B14").Value = "东方财富"
$ws.Range("C14").Value = "华建集团"
$ws.Range("A15").Value = "明阳智能"
$ws.Range("B15").Value = "立讯精密"
$ws.Range("C15").Value = "多氟多"
$ws.Range("A16").Value = "赛力斯"
$ws.Range("B16").Value = "天际股份"
$ws.Range("C16").Value = "先导智能"
$ws.Range("A17").Value = "中际旭创"
$ws.Range("B17").Value = "贵州茅台"
$ws.Range("C17").Value = "张江高科"
$ws.Range("A18").Value = "天际股份"
$ws.Range("B18").Value = "金风科技"
$ws.Range("C18").Value = "天际股份"
$ws.Range("A19").Value = "养元饮品"
$ws.Range("B19").Value = "蓝黛科技"
$ws.Range("C19").Value = "蓝丰生化"
$ws.Range("A20").Value = "长川科技"
$ws.Range("B20").Value = "养元饮品"
$ws.Range("A21").Value = "凯美特气"
$ws.Range("B21").Value = "张江高科"
$ws.Range("C21").Value = "赛力斯"
